# Module language usage added.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "R"
$ws.Range("G3").Value = "Jupyter"
$ws.Range("G4").Value = "Jupyter"
$ws.Range("G5").Value = "Java"

$ws.Range("G6").Select()
